$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This shared string is used by:
#      Overview!B2, Overview!C2, Overview!B3, Overview!C3
#      zh-cn!B2, zh-cn!B3
#      de-de!B2, de-de!B3
# ---------------------------------------------------------------------------
$newStatus = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B2").Value = $newStatus
$wsZh.Range("B3").Value = $newStatus

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B2").Value = $newStatus
$wsDe.Range("B3").Value = $newStatus

# ---------------------------------------------------------------------------
# 2. zh-cn sheet: rows 2 & 3 get "Latest Target File" (E) and
#    "Latest Handback File" (F) populated (duplicates of the source
#    file / handoff file, each with the same hyperlink target as the
#    corresponding A / C cell), and "Latest Handback DateTime" (G) gets
#    a real timestamp instead of the epoch placeholder.
# ---------------------------------------------------------------------------
$mdTarget2 = "https://github.com/OpenLocalizationTest/oltest/blob/ea9b8229e4f8f30066812a5777a86aa0be4fe34c/e2e/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md"
$xlfTargetZh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ef777b5f9c49a32dfbea3c35718b6f0c01350e61/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf"
$mdTarget3 = "https://github.com/OpenLocalizationTest/oltest/blob/ea9b8229e4f8f30066812a5777a86aa0be4fe34c/e2e/ffff3a438823-8e45-4206-ad5a-b92aae305d78.md"

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $mdTarget2, [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $xlfTargetZh, [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf")
$wsZh.Range("G2").Value = "2016-01-22 02:58:33"

$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $mdTarget3, [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md")
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $xlfTargetZh, [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.zh-cn.xlf")
$wsZh.Range("G3").Value = "2016-01-22 02:58:33"

# ---------------------------------------------------------------------------
# 3. de-de sheet: same pattern, with its own handoff-file hyperlink target
#    and its own handback timestamp.
# ---------------------------------------------------------------------------
$xlfTargetDe = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c71040fbf23ced8ac95e4d6232df08703a279021/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf"

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $mdTarget2, [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $xlfTargetDe, [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf")
$wsDe.Range("G2").Value = "2016-01-22 02:58:57"

$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $mdTarget3, [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.md")
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $xlfTargetDe, [Type]::Missing, [Type]::Missing, "ca3c0fb3-b850-4ae1-9914-d36efa6e197e.c6a2a43d5404e307aba0548c33fcaa8e348f522d.de-de.xlf")
$wsDe.Range("G3").Value = "2016-01-22 02:58:57"
